$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9660.60250133672
$ws.Range("C2").Value = 20899.18713995252
$ws.Range("D2").Value = 43012.21701260158
$ws.Range("E2").Value = 67356.68017786511

$ws.Range("B3").Value = 87187.40227170222
$ws.Range("C3").Value = 164939.3622348059
$ws.Range("D3").Value = 213192.1293490655
$ws.Range("E3").Value = 246141.7882807517

$ws.Range("B4").Value = 9416.270828340033
$ws.Range("C4").Value = 18326.9293630279
$ws.Range("D4").Value = 31818.6464115022
$ws.Range("E4").Value = 43476.87719194888

$ws.Range("B6").Value = 17406.16513793905
$ws.Range("C6").Value = 23338.21032166498
$ws.Range("D6").Value = 24646.40781747946
$ws.Range("E6").Value = 22127.86894278204

$ws.Range("B7").Value = 2678.73864901274
$ws.Range("C7").Value = 5721.422959108772
$ws.Range("D7").Value = 7349.267618257169
$ws.Range("E7").Value = 8905.534217615552

$ws.Range("B9").Value = 228125.0293063644
$ws.Range("C9").Value = 380519.3329503042
$ws.Range("D9").Value = 563468.8419364538
$ws.Range("E9").Value = 736949.3661745077

$ws.Range("B12").Value = 94860.92602788289
$ws.Range("C12").Value = 142701.9058444258
$ws.Range("D12").Value = 209030.9729086272
$ws.Range("E12").Value = 218875.257823178
